# added required experimental boolean element to valuesets
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 = "Experimental" property; it previously had no Value (column B empty).
# The ValueSet now carries the required "experimental" element, rendered as
# the literal text "true" (matching the sibling "Immutable" row's textual
# boolean rendering). A bare Value="true" assignment gets auto-coerced to a
# real Excel boolean by this engine, so force text entry with a leading
# apostrophe and then restore B7's original cell formatting (border/fill/
# wrap, taken from another cell that already carries that same style).
$ws.Range("B7").Value = "'true"
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 8 = "Date" property; the generated report timestamp moved forward.
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
